$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Duplicate the "ifoCAST long_term" sheet (the template for the new
#    "ifoCAST full" series) and place the copy right after it.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("ifoCAST long_term")
$template.Copy([System.Reflection.Missing]::Value, $template)
$new = $wb.Worksheets.Item($template.Index + 1)
$new.Name = "ifoCAST full"

# ---------------------------------------------------------------------
# 2. Update the two section headers that got new labels/text
#    ("AVERAGE1" -> "ifoCAST full", "ifo/ifoCAST" -> "ifo judgemental").
#    Re-typing the header text in Excel drops the cell's custom style,
#    so we clear formatting on the two header cells and on the other
#    two header cells that lost formatting in the same edit session.
# ---------------------------------------------------------------------
$new.Range("B2").Value = "ifoCAST full"

$new.Range("B6").Value = "ifo judgemental"
$new.Range("B6").ClearFormats()
$new.Range("B10").ClearFormats()
$new.Range("B18").ClearFormats()

# B11 ("Q0" row under "ifo judgemental") picked up the plain/top-aligned
# style used elsewhere on the sheet (style of C2) instead of its old
# center-aligned style.
$new.Range("C2").Copy()
$new.Range("B11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3. Refresh the data rows with the new "full sample" figures (N=46)
#    (written as plain decimals - this PowerShell parser chokes on
#    scientific-notation literals such as 2.45E-2)
# ---------------------------------------------------------------------
$new.Range("C3:H3").Value = @(@(0.155486188, 0.604724878, 0.919852239, 0.959089276, 0.956859551, 46))

$new.Range("C7:H7").Value = @(@(0.02458820482611912, 0.3766467348305488, 0.3986471647194397, 0.6313851160103789, 0.6378777121662803, 46))

$new.Range("C11:H11").Value = @(@(-0.05629308736242639, 1.12864044670977, 6.004745060421825, 2.450458132762489, 2.477493815665918, 46))

$new.Range("C15:H15").Value = @(@(0.035148317, 1.423122695, 9.937334559, 3.152353812, 3.18777659, 46))

$new.Range("C19:H19").Value = @(@(0.08661448288610207, 0.9270643824389406, 4.428799052064384, 2.104471204855126, 2.126447978618142, 46))

# ---------------------------------------------------------------------
# 4. Sheet "ifoCAST - ifo 2022 sample" is no longer the selected tab;
#    select the new sheet instead (cursor resting on J12, as in the
#    authored workbook) and reset its zoom back to 100%.
# ---------------------------------------------------------------------
$new.Activate()
$excel.ActiveWindow.Zoom = 100
$new.Range("J12").Select() | Out-Null

# The previously-active "ifoCAST long_term" sheet also had its
# selection extended to B1:N20 during the edit.
$template.Range("B1:N20").Select() | Out-Null
$template.Activate()

# Leave the new sheet as the final active tab.
$new.Activate()
